# ABA hold removal logic implementation
#
# Rotates the ABA-hold test-account rows in ParentCredentials /
# StudentCredentials to the next batch of generated credentials /
# ATM ticket numbers.

$wb = $excel.ActiveWorkbook

$parent  = $wb.Worksheets.Item("ParentCredentials")
$student = $wb.Worksheets.Item("StudentCredentials")

# --- Row 3 (ParentCredentials + StudentCredentials) ---------------------
$parent.Range("A3").Value = "rcg+27584685@pcci.edu"
$parent.Range("B3").Value = "rcg27584685"
$parent.Range("D3").Value = "27584685"

$student.Range("A3").Value = "ATM202111251820"
$student.Range("B3").Value = "Password@123"
$student.Range("C3").Value = "ATM202111251820 ATM202111251820"

# --- Row 4 (ParentCredentials + StudentCredentials) ---------------------
$parent.Range("A4").Value = "rcg+27584689@pcci.edu"
$parent.Range("B4").Value = "rcg27584689"
$parent.Range("D4").Value = "27584689"

$student.Range("A4").Value = "ATM202111251944"
$student.Range("C4").Value = "ATM202111251944 ATM202111251944"

# --- Row 5 (ParentCredentials + StudentCredentials) ---------------------
$parent.Range("A5").Value = "rcg+27584690@pcci.edu"
$parent.Range("B5").Value = "rcg27584690"
# D5 has no explicit text number-format in this sheet (unlike D3/D4), so
# force text storage (no leading-zero / scientific mangling) and then
# drop back to the default style so no new cell format is introduced.
$parent.Range("D5").NumberFormat = "@"
$parent.Range("D5").Value = "27584690"
$parent.Range("D5").ClearFormats()

$student.Range("A5").Value = "ATM202111251958"
$student.Range("C5").Value = "ATM202111251958 ATM202111251958"

# --- Row 6 (ParentCredentials only) --------------------------------------
$parent.Range("A6").Value = "rcg+27584691@pcci.edu"
$parent.Range("B6").Value = "rcg27584691"
$parent.Range("D6").NumberFormat = "@"
$parent.Range("D6").Value = "27584691"
$parent.Range("D6").ClearFormats()

# ParentCredentials column A was auto best-fit-sized by Excel after the
# text changed length.
$parent.Columns("A").ColumnWidth = 23.75
